$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.083.96"
$ws.Range("E2").Value = "  +0.03%  "
$ws.Range("D3").Value = "2.467.47"
$ws.Range("E3").Value = "  +0.00%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'581.69"
$ws.Range("E5").Value = "  -0.16%  "
$ws.Range("D6").Value = "'174.00"
$ws.Range("E6").Value = "  +2.67%  "
$ws.Range("E7").Value = "  -0.03%  "
$ws.Range("E8").Value = "  -0.47%  "
$ws.Range("D10").Value = "'0.165"
$ws.Range("E10").Value = "  +0.20%  "
$ws.Range("E11").Value = "  -0.41%  "
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("E13").Value = "  +0.22%  "
$ws.Range("E14").Value = "  -0.84%  "
$ws.Range("D15").Value = "66.973.79"
$ws.Range("E16").Value = "  -0.19%  "
$ws.Range("D17").Value = "2.480.55"
$ws.Range("E17").Value = "  +1.00%  "
$ws.Range("D18").Value = "'10.88"
$ws.Range("E18").Value = "  -2.15%  "
$ws.Range("D19").Value = "'7.44"
$ws.Range("E19").Value = "  -1.87%  "
$ws.Range("D20").Value = "'347.75"
$ws.Range("E20").Value = "  -1.48%  "
$ws.Range("E21").Value = "  -0.31%  "
$ws.Range("E22").Value = "  +0.00%  "
$ws.Range("D23").Value = "'69.29"
$ws.Range("E23").Value = "  +0.47%  "
$ws.Range("E24").Value = "  -1.36%  "
$ws.Range("E25").Value = "  +0.29%  "
$ws.Range("D26").Value = "'9.13"
$ws.Range("E26").Value = "  -1.13%  "
$ws.Range("D27").Value = "2.594.77"
$ws.Range("E27").Value = "  +1.07%  "
$ws.Range("D28").Value = "'1.00"
$ws.Range("E28").Value = "  +0.38%  "
$ws.Range("D29").Value = "0.0₃0898"
$ws.Range("E29").Value = "  -0.66%  "
$ws.Range("D30").Value = "'497.85"
$ws.Range("E30").Value = "  -3.91%  "
$ws.Range("D31").Value = "'7.70"
$ws.Range("E31").Value = "  -0.69%  "
$ws.Range("E32").Value = "  -0.41%  "
$ws.Range("E33").Value = "  -1.27%  "
$ws.Range("E34").Value = "  -0.04%  "
$ws.Range("E35").Value = "  +1.12%  "
$ws.Range("D36").Value = "'161.84"
$ws.Range("E36").Value = "  +1.98%  "
$ws.Range("E37").Value = "  -0.01%  "
$ws.Range("D38").Value = "'18.10"
$ws.Range("E38").Value = "  -1.46%  "
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("E40").Value = "  +0.04%  "
$ws.Range("E41").Value = "  +0.79%  "
$ws.Range("E42").Value = "  +0.03%  "
$ws.Range("D43").Value = "'4.80"
$ws.Range("E43").Value = "  +0.03%  "
$ws.Range("D44").Value = "'2.38"
$ws.Range("E44").Value = "  +0.07%  "
$ws.Range("D45").Value = "'142.15"
$ws.Range("E45").Value = "  +0.90%  "
$ws.Range("E46").Value = "  +0.48%  "
$ws.Range("E47").Value = "  -1.25%  "
$ws.Range("E48").Value = "  -1.18%  "
$ws.Range("E49").Value = "  +0.79%  "
$ws.Range("E50").Value = "  -1.69%  "
$ws.Range("E51").Value = "  +0.08%  "
